$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29, pushing existing rows 29-41 down to 30-42.
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new weekly data point.
$ws.Range("A29").Value = 3
$ws.Range("B29").Value = "Femacal de La Calera"
$ws.Range("C29").Value = "Coquimbo"
$ws.Range("D29").Value = 44529
$ws.Range("D29").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E29").Value = 5
$ws.Range("F29").Value = 100112022
$ws.Range("G29").Value = "Arveja Verde"
$ws.Range("H29").Value = "Perfection"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 73
$ws.Range("K29").Value = 17000
$ws.Range("L29").Value = 18000
$ws.Range("M29").Value = 17521
$ws.Range("N29").Value = "$/saco 25 kilos"
$ws.Range("O29").Value = "Provincia de Limarí"
$ws.Range("P29").Value = 701
$ws.Range("Q29").Value = 25
$ws.Range("R29").Value = "Hortaliza"
